$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" (F) column values to reflect the repulled/pushed data.
$ws.Range("F2").Value  = -3
$ws.Range("F7").Value  = -3
$ws.Range("F10").Value = 1
$ws.Range("F11").Value = 5
$ws.Range("F13").Value = 0
$ws.Range("F14").Value = 1
$ws.Range("F18").Value = -1
$ws.Range("F22").Value = 0
$ws.Range("F25").Value = -2
$ws.Range("F26").Value = 0
$ws.Range("F30").Value = -2
$ws.Range("F34").Value = 0
$ws.Range("F36").Value = -6
$ws.Range("F39").Value = -7
